$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventory for check stock")

$ws.Range("F6").Value = "01A012"
$ws.Range("F3").Value = "01A013"
$ws.Range("F4").Value = "01A014"
$ws.Range("F5").Value = "01A015"

$ws.Range("F2").Select()
